$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A57").Value = "GRT-USD"
